$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new lines ("line7", "line8") are inserted conceptually between "line6"
# and the "extr*" rows, pushing extr1..extr8 down by two rows; two brand new
# "extr7"/"extr8" rows are appended at the bottom with fresh data. Rewrite
# rows 8-17 directly with literal values (existing rows 8-15 get overwritten
# in place, rows 16-17 are brand new).

# Row 8: line7 (was extr1's row)
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $true

# Row 9: line8 (was extr2's row)
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = 9
$ws.Cells.Item(9, 5).Value = $true

# Row 10: extr1 (was extr3's row)
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "extr1"
$ws.Cells.Item(10, 3).Value = 5
$ws.Cells.Item(10, 4).Value = 12
$ws.Cells.Item(10, 5).Value = $true

# Row 11: extr2 (was extr4's row)
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "extr2"
$ws.Cells.Item(11, 3).Value = 5
$ws.Cells.Item(11, 4).Value = 9
$ws.Cells.Item(11, 5).Value = $true

# Row 12: extr3 (was extr5's row)
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = "extr3"
$ws.Cells.Item(12, 3).Value = 10
$ws.Cells.Item(12, 4).Value = 11
$ws.Cells.Item(12, 5).Value = $false

# Row 13: extr4 (was extr6's row)
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "extr4"
$ws.Cells.Item(13, 3).Value = 7
$ws.Cells.Item(13, 4).Value = 8
$ws.Cells.Item(13, 5).Value = $false

# Row 14: extr5 (was extr7's row)
$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = "extr5"
$ws.Cells.Item(14, 3).Value = 9
$ws.Cells.Item(14, 4).Value = 11
$ws.Cells.Item(14, 5).Value = $false

# Row 15: extr6 (was extr8's row)
$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = "extr6"
$ws.Cells.Item(15, 3).Value = 7
$ws.Cells.Item(15, 4).Value = 11
$ws.Cells.Item(15, 5).Value = $true

# Row 16: extr7 (brand new row)
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "extr7"
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = 7
$ws.Cells.Item(16, 5).Value = $false

# Row 17: extr8 (brand new row)
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "extr8"
$ws.Cells.Item(17, 3).Value = 8
$ws.Cells.Item(17, 4).Value = 5
$ws.Cells.Item(17, 5).Value = $false

# The new rows 16 & 17 need the same "A" column style (bold/centered/bordered)
# used throughout the rest of the table (e.g. row 15's A cell).
$ws.Cells.Item(15, 1).Copy()
$ws.Cells.Item(16, 1).PasteSpecial(-4122)
$ws.Cells.Item(17, 1).PasteSpecial(-4122)
